$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = 320.4349975585937
$ws.Range("G4").Value = 110.1319885253907
$ws.Range("H4").Value = 34.36952560253732

$ws.Range("C5").Value = 975.9210205078124
$ws.Range("F5").Value = 430.7210083007813
$ws.Range("H5").Value = 123.7510969708795

$ws.Range("B12").Value = 26981.6931317597
$ws.Range("C12").Value = 37313.96968745709
$ws.Range("E12").Value = 37313.96968745709
$ws.Range("G12").Value = 20766.05528958555
$ws.Range("H12").Value = 125.4904684076476
